$wb = $excel.ActiveWorkbook

# --- Update the "Conversión del día" text block on sheet "Hoja1" ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$cell = $wsHoja1.Range("A1")
$text = $cell.Value2
$text = $text -replace [regex]::Escape("✅ 1000 Bs = 7.33 = 30043.96 pesos"), "✅ 1000 Bs = 7.38 = 30258.3 pesos"
$text = $text -replace [regex]::Escape("✅ 30043.96 pesos = 7.29 = 961.62 Bs"), "✅ 30258.3 pesos = 7.37 = 978.9 Bs"
$cell.Value = $text

# --- Update the rate figures on sheet "tasas" ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 135.5
$wsTasas.Range("O10").Value = 4100
$wsTasas.Range("N12").Value = 4107.99
$wsTasas.Range("O12").Value = 132.9
